$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Fix typo "third-part tools" -> "third-party tools" in the
#     "Principles" rounded-rectangle callout (last bullet paragraph) ---
$principles = $s.Shapes.Item("Rectangle à coins arrondis 87")
$tr = $principles.TextFrame.TextRange
$lastPara = $tr.Paragraphs($tr.Paragraphs().Count)

# Route the text change through an unrelated placeholder first so the
# engine's prefix/suffix diffing doesn't fragment the run: going
# old-text -> placeholder -> new-text shares no characters with either
# end, so the final assignment lands as a single run with the run's
# original formatting (rPr) intact.
$lastPara.Text = "zzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzz"
$lastPara.Text = "Possibility to integrate third-party tools"

# --- Reposition the "Graphic 4" picture ---
# Target OOXML offsets (EMU): x=3908603, y=2304298
# Shape.Left/Top are in points and are stored single-precision (float32)
# before being floored to EMU (*12700), so feed in the float32 value
# that lands exactly on the target EMU after that conversion.
$graphic = $s.Shapes.Item("Graphic 4")
$graphic.Left = 307.7640686035156
$graphic.Top = 181.44082641601562
